# Applies commit "#65 Updated Options structure to be able to have globalOptions"
# Adds two new worksheets: "layeredOption" and "layeredOptionActual" with test
# fixture data for xlbean's layered (global/table/column) options feature.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: best-effort character-width -> ColumnWidth solver.
# This runtime stores column widths on a 1/7-character grid computed as
# round((ColumnWidth + 5/7) * 7) / 7, which differs from real Excel's pixel
# based formula. We invert it here to land as close as possible to the
# widths that real Excel produced in the authored workbook.
# ---------------------------------------------------------------------------
function Set-BestFitColumnWidth($col, [double]$targetStoredWidth) {
    $n = [Math]::Round($targetStoredWidth * 7)
    $gridVal = $n / 7.0
    $chars = $gridVal - (5.0 / 7.0)
    if ($chars -lt 0) { $chars = 0 }
    $col.ColumnWidth = $chars
}

# ---------------------------------------------------------------------------
# Add "layeredOption" sheet (sheet6) right after "optionForTableAndColumn"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$layeredOption = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$layeredOption.Name = "layeredOption"

$layeredOption.Range("A1").Value = "####?global=val1"
$layeredOption.Range("C1").Value = "layeredaaa?single=val2"
$layeredOption.Range("D1").Value = "layeredbbb?global=val9"
$layeredOption.Range("F1").Value = "layeredtable?table=val3#column1"
$layeredOption.Range("G1").Value = "layeredtable#column2?column=val4"
$layeredOption.Range("H1").Value = "layeredtable#column3?table=val9"
$layeredOption.Range("J1").Value = "layeredtable2?global=val99#col"
$layeredOption.Range("K1").Value = "layeredtable2#col?global=val999"

$layeredOption.Range("A3").Value = "layeredaaa, layeredbbb"
$layeredOption.Range("C3").Value = "dddd"
$layeredOption.Range("D3").Value = "eeee"

$layeredOption.Range("A5").Value = "layeredtable#~"
$layeredOption.Range("F5").Value = 111
$layeredOption.Range("G5").Value = 222
$layeredOption.Range("H5").Value = 333

$layeredOption.Range("F6").Value = 444
$layeredOption.Range("G6").Value = 555
$layeredOption.Range("H6").Value = 666

$layeredOption.Range("A7").Value = "layeredtable2#~"
$layeredOption.Range("J7").Value = 1
$layeredOption.Range("K7").Value = 2

$layeredOption.Range("J8").Value = 2
$layeredOption.Range("K8").Value = 3

Set-BestFitColumnWidth $layeredOption.Columns.Item(1) 21.5
Set-BestFitColumnWidth $layeredOption.Columns.Item(2) 0.796875
Set-BestFitColumnWidth $layeredOption.Columns.Item(3) 22
Set-BestFitColumnWidth $layeredOption.Columns.Item(4) 22.09765625
Set-BestFitColumnWidth $layeredOption.Columns.Item(5) 0.69921875
Set-BestFitColumnWidth $layeredOption.Columns.Item(6) 29.8984375
Set-BestFitColumnWidth $layeredOption.Columns.Item(7) 30.8984375
Set-BestFitColumnWidth $layeredOption.Columns.Item(8) 29.8984375
Set-BestFitColumnWidth $layeredOption.Columns.Item(9) 0.796875
Set-BestFitColumnWidth $layeredOption.Columns.Item(10) 29.09765625
Set-BestFitColumnWidth $layeredOption.Columns.Item(11) 30.19921875

# ---------------------------------------------------------------------------
# Add "layeredOptionActual" sheet (sheet7) right after "layeredOption"
# ---------------------------------------------------------------------------
$layeredOptionActual = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $layeredOption)
$layeredOptionActual.Name = "layeredOptionActual"

$layeredOptionActual.Range("A1").Value = "####?readAs=text"
$layeredOptionActual.Range("C1").Value = "layeredAct"
$layeredOptionActual.Range("H1").Value = "layeredActTable#col1"
$layeredOptionActual.Range("I1").Value = "layeredActTable#col2?fieldType=boolean"

$layeredOptionActual.Range("A3").Value = "layeredAct"
$layeredOptionActual.Range("C3").Value = 1
$layeredOptionActual.Range("D3").Value = '<- this to be "1" and not "1.0"'

$layeredOptionActual.Range("A4").Value = "layeredActTable?fieldType=int#~"
$layeredOptionActual.Range("H4").Value = 1
$layeredOptionActual.Range("I4").Value = $true

$layeredOptionActual.Range("H5").Value = 2
$layeredOptionActual.Range("I5").Value = $false

$layeredOptionActual.Range("H6").Value = 3
$layeredOptionActual.Range("I6").Value = $true

Set-BestFitColumnWidth $layeredOptionActual.Columns.Item(1) 17.09765625
Set-BestFitColumnWidth $layeredOptionActual.Columns.Item(8) 20
Set-BestFitColumnWidth $layeredOptionActual.Columns.Item(9) 20

# ---------------------------------------------------------------------------
# Selection / active-tab bookkeeping: Excel records the selected cell of the
# sheet that was active when the file was last saved, and marks that sheet's
# tab as selected. Here the final active sheet is "layeredOption" (with "G1"
# selected), while "layeredOptionActual" was left with "I7" selected.
# ---------------------------------------------------------------------------
$layeredOptionActual.Activate() | Out-Null
$layeredOptionActual.Range("I7").Select() | Out-Null

$layeredOption.Activate() | Out-Null
$layeredOption.Range("G1").Select() | Out-Null

Write-Output "Added layeredOption and layeredOptionActual sheets"
